$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1457.1111
$ws.Range("H131").Value = 2950.8235
$ws.Range("I131").Value = 257.375
$ws.Range("J131").Value = 5345
$ws.Range("K131").Value = 772.125
$ws.Range("L131").Value = 16035
$ws.Range("M131").Value = 4267.875
$ws.Range("N131").Value = -26115
$ws.Range("H135").Value = 1127.3125
$ws.Range("I135").Value = 625.1111
$ws.Range("J135").Value = 3839.2
$ws.Range("K135").Value = 5625.9999
$ws.Range("L135").Value = 34552.8
$ws.Range("M135").Value = -3090.9999
$ws.Range("N135").Value = -39622.8
$ws.Range("H137").Value = 1557.641
$ws.Range("I137").Value = 1436.7916
$ws.Range("J137").Value = 1751
$ws.Range("K137").Value = 4310.3748
$ws.Range("L137").Value = 5253
$ws.Range("M137").Value = -1760.3748
$ws.Range("N137").Value = -10353
$ws.Range("H138").Value = 2291.0344
$ws.Range("I138").Value = 1821.138
$ws.Range("J138").Value = 2525.9827
$ws.Range("K138").Value = 5463.414
$ws.Range("L138").Value = 7577.9481
$ws.Range("M138").Value = -323.4139999999998
$ws.Range("N138").Value = -17857.9481

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 963.7143
$ws.Range("I74").Value = 710.3125
$ws.Range("J74").Value = 3666.6667
$ws.Range("K74").Value = 710.3125
$ws.Range("L74").Value = 3666.6667
$ws.Range("M74").Value = 163.6875
$ws.Range("N74").Value = -5414.6667
$ws.Range("H77").Value = 963.7143
$ws.Range("I77").Value = 710.3125
$ws.Range("J77").Value = 3666.6667
$ws.Range("K77").Value = 3551.5625
$ws.Range("L77").Value = 18333.3335
$ws.Range("M77").Value = 816.4375
$ws.Range("N77").Value = -27069.3335

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 19987
$ws.Range("J35").Value = 19987
$ws.Range("L35").Value = 19987
$ws.Range("N35").Value = -20607
$ws.Range("H82").Value = 26996.309
$ws.Range("J82").Value = 43542.145
$ws.Range("L82").Value = 43542.145
$ws.Range("N82").Value = -44308.145
$ws.Range("H85").Value = 26996.309
$ws.Range("J85").Value = 43542.145
$ws.Range("L85").Value = 43542.145
$ws.Range("N85").Value = -46194.145

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5146.2666
$ws.Range("I31").Value = 1365.8948
$ws.Range("K31").Value = 1365.8948
$ws.Range("M31").Value = -1070.8948
$ws.Range("H34").Value = 5146.2666
$ws.Range("I34").Value = 1365.8948
$ws.Range("K34").Value = 1365.8948
$ws.Range("M34").Value = -1163.8948
$ws.Range("H41").Value = 7285.4287
$ws.Range("I41").Value = 3000
$ws.Range("J41").Value = 8999.6
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 8999.6
$ws.Range("M41").Value = -2572
$ws.Range("N41").Value = -9855.6
$ws.Range("H50").Value = 18599.2
$ws.Range("J50").Value = 18599.2
$ws.Range("L50").Value = 18599.2
$ws.Range("N50").Value = -19849.2
$ws.Range("H51").Value = 16665.834
$ws.Range("J51").Value = 16665.834
$ws.Range("L51").Value = 16665.834
$ws.Range("N51").Value = -18137.834
$ws.Range("H59").Value = 30399
$ws.Range("J59").Value = 30399
$ws.Range("L59").Value = 30399
$ws.Range("N59").Value = -32689
$ws.Range("H60").Value = 10448
$ws.Range("J60").Value = 10448
$ws.Range("L60").Value = 10448
$ws.Range("N60").Value = -11470
$ws.Range("H61").Value = 16665.834
$ws.Range("J61").Value = 16665.834
$ws.Range("L61").Value = 16665.834
$ws.Range("N61").Value = -17361.834
$ws.Range("H68").Value = 19408.268
$ws.Range("J68").Value = 19408.268
$ws.Range("L68").Value = 19408.268
$ws.Range("N68").Value = -20906.268
$ws.Range("H71").Value = 19408.268
$ws.Range("J71").Value = 19408.268
$ws.Range("L71").Value = 58224.804
$ws.Range("N71").Value = -65712.804

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 12500471
$ws.Range("J34").Value = 13889392
$ws.Range("L34").Value = 41668176
$ws.Range("N34").Value = -41668344
$ws.Range("H39").Value = 1474.4615
$ws.Range("J39").Value = 1516.2972
$ws.Range("L39").Value = 4548.8916
$ws.Range("N39").Value = -5136.8916
$ws.Range("H55").Value = 1442.3077
$ws.Range("J55").Value = 1495.8334
$ws.Range("L55").Value = 4487.5002
$ws.Range("N55").Value = -4841.5002
$ws.Range("H122").Value = 3024.561
$ws.Range("I122").Value = 398
$ws.Range("K122").Value = 3582
$ws.Range("M122").Value = -1132

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8761.875
$ws.Range("H46").Value = 4062.4375
$ws.Range("J46").Value = 4062.4375
$ws.Range("L46").Value = 4062.4375
$ws.Range("N46").Value = -4374.4375
$ws.Range("H80").Value = 29967060
$ws.Range("J80").Value = 201999.5
$ws.Range("L80").Value = 201999.5
$ws.Range("N80").Value = -203995.5
$ws.Range("H83").Value = 29967060
$ws.Range("J83").Value = 201999.5
$ws.Range("L83").Value = 1009997.5
$ws.Range("N83").Value = -1019981.5
$ws.Range("H126").Value = 1754.6666
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 1825.6
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 5476.799999999999
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -10416.8
$ws.Range("H132").Value = 1971.7858
$ws.Range("I132").Value = 1215.9166
$ws.Range("J132").Value = 6507
$ws.Range("K132").Value = 3647.7498
$ws.Range("L132").Value = 19521
$ws.Range("M132").Value = -1117.7498
$ws.Range("N132").Value = -24581

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1521.6774
$ws.Range("I68").Value = 1464.3043
$ws.Range("K68").Value = 1464.3043
$ws.Range("M68").Value = -715.3043
$ws.Range("H71").Value = 1521.6774
$ws.Range("I71").Value = 1464.3043
$ws.Range("K71").Value = 7321.5215
$ws.Range("M71").Value = -3577.5215

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25496
$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26716
$ws.Range("H81").Value = 3362.4707
$ws.Range("J81").Value = 2980.2
$ws.Range("L81").Value = 5960.4
$ws.Range("N81").Value = -8082.4
$ws.Range("H84").Value = 3362.4707
$ws.Range("J84").Value = 2980.2
$ws.Range("L84").Value = 29802
$ws.Range("N84").Value = -40410
$ws.Range("H100").Value = 2286.111
$ws.Range("I100").Value = 2465.3333
$ws.Range("J100").Value = 1927.6666
$ws.Range("K100").Value = 4930.6666
$ws.Range("L100").Value = 3855.3332
$ws.Range("M100").Value = -4389.6666
$ws.Range("N100").Value = -4937.3332
$ws.Range("H107").Value = 621.8378
$ws.Range("I107").Value = 625.0345
$ws.Range("J107").Value = 610.25
$ws.Range("K107").Value = 1875.1035
$ws.Range("L107").Value = 1830.75
$ws.Range("M107").Value = 44.89650000000006
$ws.Range("N107").Value = -5670.75
$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1130
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 1836.2874
$ws.Range("I136").Value = 1527.6104
$ws.Range("J136").Value = 4213.1
$ws.Range("K136").Value = 4582.831200000001
$ws.Range("L136").Value = 12639.3
$ws.Range("M136").Value = -2032.831200000001
$ws.Range("N136").Value = -17739.3
